$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 6-9 (MuSCs/Resolving-Mac -> ECs/Resolving-Mac pairs)
$ws.Rows("6:9").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 15.428109
$ws.Range("H2").Value = 46.284327
$ws.Range("I2").Value = 0.105145687357564
$ws.Range("J2").Value = 0.105145687357564
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1790523333333333
$ws.Range("N2").Value = 0.537157
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2.762438915371
$ws.Range("R2").Value = 24.861950238339
$ws.Range("S2").Value = 0.105145687357564
$ws.Range("T2").Value = 0.105145687357564

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 72.11798333333333
$ws.Range("H3").Value = 216.35395
$ws.Range("I3").Value = 0.4914986618531588
$ws.Range("J3").Value = 0.4914986618531588
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1790523333333333
$ws.Range("N3").Value = 0.537157
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 12.91289319112778
$ws.Range("R3").Value = 116.21603872015
$ws.Range("S3").Value = 0.4914986618531588
$ws.Range("T3").Value = 0.4914986618531588

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 17.753286
$ws.Range("H4").Value = 53.25985799999999
$ws.Range("I4").Value = 0.120992239510715
$ws.Range("J4").Value = 0.120992239510715
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1790523333333333
$ws.Range("N4").Value = 0.537157
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 3.178767282634
$ws.Range("R4").Value = 28.608905543706
$ws.Range("S4").Value = 0.120992239510715
$ws.Range("T4").Value = 0.120992239510715

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 41.43140433333334
$ws.Range("H5").Value = 124.294213
$ws.Range("I5").Value = 0.2823634112785623
$ws.Range("J5").Value = 0.2823634112785622
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1790523333333333
$ws.Range("N5").Value = 0.537157
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 7.418389619160113
$ws.Range("R5").Value = 66.76550657244101
$ws.Range("S5").Value = 0.2823634112785623
$ws.Range("T5").Value = 0.2823634112785622
